$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 218, pushing the existing rows 218-252 down to 219-253
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new weekly price-report entry
$ws.Range("A218").Value = 11
$ws.Range("B218").Value = "Vega Monumental Concepción"
$ws.Range("C218").Value = "Bíobío"
$ws.Range("D218").Value = 45244
$ws.Range("E218").Value = 8
$ws.Range("F218").Value = 100112043
$ws.Range("G218").Value = "Pepino ensalada"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 150
$ws.Range("K218").Value = 17000
$ws.Range("L218").Value = 17000
$ws.Range("M218").Value = 17000
$ws.Range("N218").Value = "$/caja 50 unidades"
$ws.Range("O218").Value = "Región de Arica y Parinacota"
$ws.Range("P218").Value = 340
$ws.Range("Q218").Value = 50
$ws.Range("R218").Value = "Hortaliza"
